$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial number that needs to move
# forward by one day (46081 -> 46082, i.e. 2026-02-28 -> 2026-03-01) for
# every data row (rows 2 through 75).
for ($row = 2; $row -le 75; $row++) {
    $ws.Cells.Item($row, 3).Value2 = 46082
}
